$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 272; this shifts old rows 272-286 down to 273-287
# (carrying all their existing data/formatting with them).
$ws.Rows.Item(272).Insert()

# Populate the newly inserted row 272 with the new weekly record.
$ws.Cells.Item(272, 1).Value = 10
$ws.Cells.Item(272, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(272, 3).Value = "La Araucanía"
$ws.Cells.Item(272, 4).Value = 44585
$ws.Cells.Item(272, 5).Value = 9
$ws.Cells.Item(272, 6).Value = "Fruta"
$ws.Cells.Item(272, 7).Value = 100108
$ws.Cells.Item(272, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(272, 9).Value = 100108002
$ws.Cells.Item(272, 10).Value = "Mango"
$ws.Cells.Item(272, 11).Value = "Sin especificar"
$ws.Cells.Item(272, 12).Value = "Primera"
$ws.Cells.Item(272, 13).Value = 1200
$ws.Cells.Item(272, 14).Value = 8000
$ws.Cells.Item(272, 15).Value = 8000
$ws.Cells.Item(272, 16).Value = 8000
$ws.Cells.Item(272, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(272, 18).Value = "Perú"
$ws.Cells.Item(272, 19).Value = 2000
$ws.Cells.Item(272, 20).Value = 4
